# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Macroferia Regional de Talca" / Naranja
# right before the current row 414, shifting the existing rows 414-440 down
# to 415-441 (dimension grows from A1:T440 to A1:T441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 414 (pushes 414..440 down to 415..441).
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with this week's reading.
$ws.Range("A414").Value = 5
$ws.Range("B414").Value = "Macroferia Regional de Talca"
$ws.Range("C414").Value = "Maule"
$ws.Range("D414").Value = 44610
$ws.Range("E414").Value = 7
$ws.Range("F414").Value = "Fruta"
$ws.Range("G414").Value = 100102
$ws.Range("H414").Value = "Cítricos"
$ws.Range("I414").Value = 100102005
$ws.Range("J414").Value = "Naranja"
$ws.Range("K414").Value = "Valencia"
$ws.Range("L414").Value = "Primera"
$ws.Range("M414").Value = 300
$ws.Range("N414").Value = 9000
$ws.Range("O414").Value = 9000
$ws.Range("P414").Value = 9000
$ws.Range("Q414").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R414").Value = "Región de O'Higgins"
$ws.Range("S414").Value = 600
$ws.Range("T414").Value = 15
